$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Correct the freight forwarder company name in cell E2
# (東風航空貨運承攬有限公司 -> 東方航空貨運承攬有限公司)
$ws.Range("E2").Value = "東方航空貨運承攬有限公司"
$ws.Range("E2").Select()
